$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 215.44444
$ws.Range("I33").Value = 215.44444
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 215.44444
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 13.55556000000001
$ws.Range("H57").Value = 40000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 40000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 120000
$ws.Range("N57").Value = -120998
$ws.Range("H69").Value = 13166.667
$ws.Range("I69").Value = 4200
$ws.Range("J69").Value = 14960
$ws.Range("K69").Value = 12600
$ws.Range("L69").Value = 44880
$ws.Range("M69").Value = -11726
$ws.Range("N69").Value = -46628
$ws.Range("H72").Value = 13166.667
$ws.Range("I72").Value = 4200
$ws.Range("J72").Value = 14960
$ws.Range("K72").Value = 37800
$ws.Range("L72").Value = 134640
$ws.Range("M72").Value = -33432
$ws.Range("N72").Value = -143376
$ws.Range("H80").Value = 843.2963
$ws.Range("J80").Value = 741.6923
$ws.Range("L80").Value = 2225.0769
$ws.Range("N80").Value = -4221.0769
$ws.Range("H83").Value = 843.2963
$ws.Range("J83").Value = 741.6923
$ws.Range("L83").Value = 6675.2307
$ws.Range("N83").Value = -16659.2307
$ws.Range("H112").Value = 3150
$ws.Range("J112").Value = 3150
$ws.Range("L112").Value = 9450
$ws.Range("N112").Value = -11666
$ws.Range("H113").Value = 2141.2856
$ws.Range("I113").Value = 1997.25
$ws.Range("K113").Value = 1997.25
$ws.Range("M113").Value = 1256.75
$ws.Range("H137").Value = 2392.9614
$ws.Range("I137").Value = 2319.5264
$ws.Range("K137").Value = 6958.5792
$ws.Range("M137").Value = -4408.5792
$ws.Range("H138").Value = 4787.294
$ws.Range("I138").Value = 3291.5715
$ws.Range("J138").Value = 5175.074
$ws.Range("K138").Value = 9874.7145
$ws.Range("L138").Value = 15525.222
$ws.Range("M138").Value = -4734.7145
$ws.Range("N138").Value = -25805.222
$ws.Range("N33").ClearContents()
$ws.Range("M57").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2415.3333
$ws.Range("J45").Value = 2000
$ws.Range("L45").Value = 2000
$ws.Range("N45").Value = -2754
$ws.Range("H61").Value = 3277.0908
$ws.Range("I61").Value = 3354.8
$ws.Range("K61").Value = 3354.8
$ws.Range("M61").Value = -3142.8
$ws.Range("H74").Value = 1469.0555
$ws.Range("I74").Value = 1406.6428
$ws.Range("K74").Value = 1406.6428
$ws.Range("M74").Value = -532.6428000000001
$ws.Range("H77").Value = 1469.0555
$ws.Range("I77").Value = 1406.6428
$ws.Range("K77").Value = 7033.214
$ws.Range("M77").Value = -2665.214
$ws.Range("H122").Value = 3787.2666
$ws.Range("I122").Value = 2103.6
$ws.Range("J122").Value = 5470.933
$ws.Range("K122").Value = 6310.799999999999
$ws.Range("L122").Value = 16412.799
$ws.Range("M122").Value = -3860.799999999999
$ws.Range("N122").Value = -21312.799
$ws.Range("H130").Value = 36380
$ws.Range("J130").Value = 36380
$ws.Range("L130").Value = 36380
$ws.Range("N130").Value = -46420
$ws.Range("H136").Value = 3277.0908
$ws.Range("I136").Value = 3354.8
$ws.Range("K136").Value = 10064.4
$ws.Range("M136").Value = -7514.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3565.2903
$ws.Range("I20").Value = 3317.45
$ws.Range("J20").Value = 4015.9092
$ws.Range("K20").Value = 3317.45
$ws.Range("L20").Value = 4015.9092
$ws.Range("M20").Value = -3070.45
$ws.Range("N20").Value = -4509.9092
$ws.Range("H22").Value = 275
$ws.Range("I22").Value = 275
$ws.Range("K22").Value = 275
$ws.Range("M22").Value = -102
$ws.Range("H60").Value = 66437.25
$ws.Range("J60").Value = 66437.25
$ws.Range("L60").Value = 66437.25
$ws.Range("N60").Value = -67635.25
$ws.Range("H75").Value = 49999
$ws.Range("J75").Value = 49999
$ws.Range("L75").Value = 49999
$ws.Range("N75").Value = -51871
$ws.Range("H78").Value = 49999
$ws.Range("J78").Value = 49999
$ws.Range("L78").Value = 149997
$ws.Range("N78").Value = -159357
$ws.Range("H139").Value = 89993.336
$ws.Range("J139").Value = 89993.336
$ws.Range("L139").Value = 89993.336
$ws.Range("N139").Value = -100273.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 114169
$ws.Range("I31").Value = 1420.3334
$ws.Range("K31").Value = 1420.3334
$ws.Range("M31").Value = -1125.3334
$ws.Range("H34").Value = 114169
$ws.Range("I34").Value = 1420.3334
$ws.Range("K34").Value = 1420.3334
$ws.Range("M34").Value = -1218.3334
$ws.Range("H52").Value = 74569.336
$ws.Range("I52").Value = 77709
$ws.Range("J52").Value = 72999.5
$ws.Range("K52").Value = 77709
$ws.Range("L52").Value = 72999.5
$ws.Range("M52").Value = -77415
$ws.Range("N52").Value = -73587.5
$ws.Range("H93").Value = 8599.6
$ws.Range("I93").Value = 8599.6
$ws.Range("K93").Value = 8599.6
$ws.Range("M93").Value = -6727.6
$ws.Range("H139").Value = 89995
$ws.Range("J139").Value = 89995
$ws.Range("L139").Value = 89995
$ws.Range("N139").Value = -100275

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 8305.8125
$ws.Range("J39").Value = 10741.5
$ws.Range("L39").Value = 32224.5
$ws.Range("N39").Value = -32812.5
$ws.Range("H44").Value = 350
$ws.Range("I44").Value = 350
$ws.Range("K44").Value = 1050
$ws.Range("M44").Value = -652
$ws.Range("H58").Value = 28209.889
$ws.Range("I58").Value = 945
$ws.Range("J58").Value = 35999.855
$ws.Range("K58").Value = 2835
$ws.Range("L58").Value = 107999.565
$ws.Range("M58").Value = -2707
$ws.Range("N58").Value = -108255.565
$ws.Range("H87").Value = 19002
$ws.Range("I87").Value = 19002
$ws.Range("K87").Value = 57006
$ws.Range("M87").Value = -55758
$ws.Range("H90").Value = 19002
$ws.Range("I90").Value = 19002
$ws.Range("K90").Value = 171018
$ws.Range("M90").Value = -164778
$ws.Range("H107").Value = 88574.586
$ws.Range("I107").Value = 1184.6
$ws.Range("J107").Value = 150996
$ws.Range("K107").Value = 3553.8
$ws.Range("L107").Value = 452988
$ws.Range("M107").Value = -1633.8
$ws.Range("N107").Value = -456828

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7150.276
$ws.Range("I70").Value = 5821.952
$ws.Range("K70").Value = 5821.952
$ws.Range("M70").Value = -5551.952
$ws.Range("H73").Value = 7150.276
$ws.Range("I73").Value = 5821.952
$ws.Range("K73").Value = 5821.952
$ws.Range("M73").Value = -4885.952
$ws.Range("H80").Value = 836881.8
$ws.Range("I80").Value = 628934.1
$ws.Range("J80").Value = 1252777.1
$ws.Range("K80").Value = 628934.1
$ws.Range("L80").Value = 1252777.1
$ws.Range("M80").Value = -627936.1
$ws.Range("N80").Value = -1254773.1
$ws.Range("H83").Value = 836881.8
$ws.Range("I83").Value = 628934.1
$ws.Range("J83").Value = 1252777.1
$ws.Range("K83").Value = 3144670.5
$ws.Range("L83").Value = 6263885.5
$ws.Range("M83").Value = -3139678.5
$ws.Range("N83").Value = -6273869.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("H128").Value = 100429
$ws.Range("J128").Value = 100429
$ws.Range("L128").Value = 100429
$ws.Range("N128").Value = -110389
$ws.Range("H132").Value = 6516.125
$ws.Range("I132").Value = 5518.4287
$ws.Range("K132").Value = 16555.2861
$ws.Range("M132").Value = -14025.2861
$ws.Range("N25").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 98489.664
$ws.Range("J124").Value = 98489.664
$ws.Range("L124").Value = 98489.664
$ws.Range("N124").Value = -108309.664
$ws.Range("H132").Value = 16780.395
$ws.Range("I132").Value = 2104.0908
$ws.Range("J132").Value = 67230.19
$ws.Range("K132").Value = 6312.2724
$ws.Range("L132").Value = 201690.57
$ws.Range("M132").Value = -3782.2724
$ws.Range("N132").Value = -206750.57
$ws.Range("H138").Value = 125000
$ws.Range("J138").Value = 125000
$ws.Range("L138").Value = 125000
$ws.Range("N138").Value = -135280
